$d = $word.ActiveDocument

$pairs = @(
    @("613×6=", "914×9="),
    @("922×6=", "509×7="),
    @("350×5=", "281×4="),
    @("385×9=", "822×2="),
    @("851×8=", "338×7="),
    @("342×8=", "112×4="),
    @("531×8=", "352×7="),
    @("757×7=", "207×5="),
    @("492×6=", "631×7="),
    @("538×6=", "696×2="),
    @("118×2=", "840×2="),
    @("450×2=", "133×7="),
    @("284×3=", "641×8="),
    @("721×6=", "131×4="),
    @("390×9=", "385×2="),
    @("192×8=", "990×4="),
    @("316×6=", "777×9="),
    @("351×6=", "567×4="),
    @("998×4=", "755×5="),
    @("105×9=", "675×9="),
    @("494×4=", "730×3="),
    @("774×6=", "472×2="),
    @("814×7=", "961×5="),
    @("372×5=", "124×7="),
    @("186×8=", "686×6=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
